$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns being updated, in order, mapped to column letters B, D, E, F, G, H, I, K, M
$cols = @("B","D","E","F","G","H","I","K","M")

# New loading_percent values for rows 2..25 (corresponding to A column index 0..23),
# one row per array, values aligned with $cols above.
$data = @(
    @(5.457444107075544,4.846976662484721,11.95299178546759,26.70639490867826,33.43590813472551,15.04613355717483,20.82625208142185,10.31713435585459,14.16767760823589),
    @(5.361873318003035,4.846677358962811,11.81774575190637,26.58457021936836,33.16175406992028,15.06726363783326,20.92530030009729,9.924743154793708,13.92647148850223),
    @(5.303175708217762,4.846566277226334,11.73813056351785,26.51841433758723,33.00566386869187,15.08372687635955,20.98982966026263,9.677773283179794,13.78008689841568),
    @(5.279283419172097,4.846539598065029,11.70658432574024,26.49364764900585,32.94519627456081,15.09131063402525,21.01705954635498,9.575774891639732,13.72094713737711),
    @(5.275318716991836,4.846536299819838,11.70140131705535,26.48966804269024,32.93534702141622,15.09262266988015,21.02163743605476,9.558761309571203,13.71116039071724),
    @(5.302853335016533,4.846565841768847,11.73770143970351,26.5180714268746,33.00483558958624,15.08382561496304,20.99019311186836,9.676402961025403,13.77928713387838),
    @(5.42451361055666,4.84685855790824,11.90567006693613,26.66260770538665,33.33887867581507,15.05269355584179,20.85963313249632,10.18318393262271,14.08420226509173),
    @(5.661596384719958,4.847996277709884,12.26055941586339,27.01365237542366,34.0879788232872,15.01943241333521,20.63307343190866,11.25317488195368,14.69219410118009),
    @(5.833110687167629,4.849158946481009,12.53447292127511,27.3112266814294,34.69108125656733,15.01205977917454,20.48458984023916,12.06402323933676,15.14030167380972),
    @(5.910206015043779,4.849755414819076,12.66141708915778,27.45481407778755,34.9758089825937,15.01243090814978,20.42094502830577,12.41356527500783,15.34350165424188),
    @(5.939239130231468,4.849990720076963,12.70977904891434,27.51033079744132,35.08502209329519,15.01310808928406,20.39740581545081,12.54313082747676,15.42027845957346),
    @(5.932993927388623,4.849939627599282,12.69935121902338,27.4983240732828,35.0614407509725,15.01293836691131,20.40245041237126,12.51535129901675,15.40375199086306),
    @(5.912597972213937,4.849774585998452,12.66539028590681,27.45935878910049,34.98476667173846,15.01247586073148,20.41899718152189,12.42428089441244,15.34982198090168),
    @(5.900083064243355,4.849674713867023,12.64462477537099,27.43563912593424,34.93797999688449,15.01226247235864,20.42920572727088,12.36813268543176,15.31676382366329),
    @(5.828051250718607,4.849121301326988,12.52622001419772,27.30200505889907,34.67267399813576,15.01211057449713,20.48882773896465,12.04078924107483,15.12700235785209),
    @(5.783604266240537,4.848798903889779,12.45414943025948,27.22210341724674,34.51250500101747,15.01297226985641,20.52640357446191,11.83501055030108,15.01036802523038),
    @(5.757953823488768,4.848619845085196,12.41291928357946,27.17692224108613,34.42136256957573,15.01381846991146,20.5483833458235,11.71483472321337,14.94322665584784),
    @(5.749255175312677,4.848560322559059,12.3989991393752,27.16175910689969,34.39067501720704,15.01416515419138,20.55588835532896,11.6738338428513,14.92048647596733),
    @(5.788344800944692,4.848832565821417,12.46179872193425,27.23052902782926,34.52945428049777,15.01284424922776,20.52236555947838,11.85710430153924,15.02279035759678),
    @(5.918593348684134,4.849822808829331,12.67535790104685,27.47077312991393,35.00725070677386,15.01259713955645,20.41412174399528,12.45110656568644,15.36566776538214),
    @(6.002766755029108,4.850524863487889,12.81660607546281,27.63443080270725,35.32758940590207,15.01556380563015,20.34665252770256,12.82300412750557,15.58872983230327),
    @(5.957937714093182,4.850145233348015,12.74108101900977,27.54648909662397,35.15591348244369,15.01369397860569,20.38236225836569,12.62601386670993,15.46979628341103),
    @(5.78620190609518,4.848817327604422,12.45833983770299,27.22671745341567,34.52178857291781,15.01290103462568,20.52418997194269,11.84712154260891,15.01717448534512),
    @(5.597804806368887,4.847629917476141,12.16206187256214,26.91159474270883,33.87571074696428,15.02544140184535,20.69120874989537,10.93761608353851,14.5271325536261)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $rowValues[$j]
    }
}
